$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.229.34'
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').Value = '3.340.49'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''555.88'
$ws.Range('E5').Value = '  +0.17%  '
$ws.Range('D6').Value = '''174.11'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('E7').Value = '  +1.05%  '
$ws.Range('D8').Value = '3.334.03'
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  +0.06%  '
$ws.Range('D10').Value = '''0.170'
$ws.Range('E10').Value = '  +6.77%  '
$ws.Range('E11').Value = '  +2.19%  '
$ws.Range('D12').Value = '''53.97'
$ws.Range('E12').Value = '  -0.21%  '
$ws.Range('D13').Value = '''0.0000278'
$ws.Range('E13').Value = '  +2.77%  '
$ws.Range('E14').Value = '  +0.78%  '
$ws.Range('D15').Value = '3.879.06'
$ws.Range('E15').Value = '  +0.20%  '
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('D17').Value = '''18.17'
$ws.Range('E17').Value = '  -0.65%  '
$ws.Range('D18').Value = '3.341.81'
$ws.Range('E18').Value = '  +0.10%  '
$ws.Range('D19').Value = '65.082.13'
$ws.Range('E19').Value = '  +2.05%  '
$ws.Range('D20').Value = '''11.76'
$ws.Range('E20').Value = '  -0.66%  '
$ws.Range('E21').Value = '  +1.46%  '
$ws.Range('D22').Value = '''454.29'
$ws.Range('E22').Value = '  +5.86%  '
$ws.Range('D23').Value = '''4.94'
$ws.Range('E23').Value = '  +7.70%  '
$ws.Range('D24').Value = '''4.09'
$ws.Range('E24').Value = '  -0.36%  '
$ws.Range('D25').Value = '''13.92'
$ws.Range('E25').Value = '  +6.69%  '
$ws.Range('D26').Value = '''87.08'
$ws.Range('E26').Value = '  +3.70%  '
$ws.Range('D27').Value = '''2.87'
$ws.Range('E27').Value = '  +2.17%  '
$ws.Range('D28').Value = '''10.75'
$ws.Range('E28').Value = '  +1.28%  '
$ws.Range('D29').Value = '''8.66'
$ws.Range('E29').Value = '  -0.59%  '
$ws.Range('E30').Value = '  +4.97%  '
$ws.Range('D31').Value = '''6.61'
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('E32').Value = '  +0.55%  '
$ws.Range('D33').Value = '''567.56'
$ws.Range('E33').Value = '  -3.83%  '
$ws.Range('D34').Value = '''60.91'
$ws.Range('E34').Value = '  +4.06%  '
$ws.Range('E35').Value = '  +0.43%  '
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = '''3.57'
$ws.Range('E37').Value = '  +3.03%  '
$ws.Range('E38').Value = '  -1.86%  '
$ws.Range('D39').Value = '''35.37'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('D40').Value = '0.0₃0737'
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').Value = '''0.368'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('D42').Value = '3.068.83'
$ws.Range('E42').Value = '  -1.28%  '
$ws.Range('E43').Value = '  -1.11%  '
$ws.Range('D44').Value = '''0.0415'
$ws.Range('E44').Value = '  +2.59%  '
$ws.Range('E45').Value = '  +0.97%  '
$ws.Range('E46').Value = '  +3.41%  '
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('D49').Value = '''141.54'
$ws.Range('E49').Value = '  +4.76%  '
$ws.Range('E50').Value = '  -2.38%  '
$ws.Range('E51').Value = '  -0.11%  '
